# The DTR sheet had placeholder "01:00:00" values left over in the
# "OFFICIAL BUSINESS DEPARTURE" (K) and "OFFICIAL BUSINESS ARRIVAL" (N)
# columns for every data row. Clear them out, leaving the
# "OFFICIAL BUSINESS TIME START"/"TIME END" (L/M) values untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EMPLOYEE DTR")

$ws.Range("K8:K15").ClearContents()
$ws.Range("N8:N15").ClearContents()
